$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert a brand-new slide ("Other Cool Things About R Markdown") right
#    before the existing "How To Learn R Markdown" slide (currently #19).
#    The simplest reliable way to get a Title+Content placeholder layout
#    identical to its neighbour is to duplicate slide 19 and then move the
#    duplicate in front of the original, so afterwards:
#       #19 -> new slide (to be rewritten below)
#       #20 -> original "How To Learn R Markdown" slide (content tweaked)
#       #21 -> "Requirements For Your Computer" (unchanged, just shifted)
# ---------------------------------------------------------------------------
$learnSlide = $p.Slides.Item(19)
$dupRange = $learnSlide.Duplicate()
$newSlide = $dupRange.Item(1)
$newSlide.MoveTo(19)

# ---------------------------------------------------------------------------
# 2. Rewrite the new slide #19 with the "Other Cool Things" content.
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Item(19)

$newTitle = $newSlide.Shapes.Item(1)
$newTitle.TextFrame.TextRange.Text = "Other Cool Things About R Markdown"

$newBody = $newSlide.Shapes.Item(2)
$newBody.Left = 66
$newBody.Top = 143.74992125984252
$newBody.Width = 828
$newBody.Height = 380.25
$newBody.TextFrame.AutoSize = 2

$newBodyTr = $newBody.TextFrame.TextRange
$newLines = @(
    "Knitting run a fresh R session, increasing reproducibility",
    "R Notebooks are even more responsive than R Markdown pure",
    "R Presentations can make presentations that are dynamic and flexible, just like R Markdown documents (but learning curve!)",
    "It is possible to run other programming languages in code chunks, like Python",
    "It is possible to create APA formatted documents (for submission even) with templates",
    "Automatic citation and bibliography generation are possible",
    "Works extremely well with version control software (VCS) such as git and github!"
)
$newBodyTr.Text = [string]::Join([char]13, $newLines)

# Split the final bullet's "git"/"github" mentions into their own runs,
# mirroring how PowerPoint isolates words flagged by the spell checker.
$lastPara = $newBodyTr.Paragraphs(7)
$gitRange = $newBodyTr.Characters($lastPara.Start + 65, 3)
$gitRange.Text = "git"
$githubRange = $newBodyTr.Characters($lastPara.Start + 65 + 3 + 5, 6)
$githubRange.Text = "github"

# ---------------------------------------------------------------------------
# 3. Tweak the original "How To Learn R Markdown" slide, now at #20:
#    - First bullet becomes three runs (middle one bold) and its ending
#      changes from "in it." to "in R Markdown."
#    - A new final sub-bullet "Drop by my office." is appended.
# ---------------------------------------------------------------------------
$learnSlide = $p.Slides.Item(20)
$learnBody = $learnSlide.Shapes.Item(2)
$learnTr = $learnBody.TextFrame.TextRange

$learnTr.InsertAfter([char]13 + "Drop by my office.")

$firstPara = $learnTr.Paragraphs(1)
$firstPara.Text = "Go try it! Commit to doing 1 full analysis in R Markdown."
$boldRange = $learnTr.Characters($firstPara.Start + 11, 32)
$boldRange.Font.Bold = $true

$firstPara = $learnTr.Paragraphs(1)
$suffixStart = $firstPara.Start + 11 + 32
$suffixLen = $firstPara.Length - 11 - 32
$suffixRange = $learnTr.Characters($suffixStart, $suffixLen)
$suffixRange.Text = "in R Markdown."
